$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").Value = 0.8238193031958776
$ws.Range("C6").Value = 0.03381894895968136
$ws.Range("D6").Value = 0.8024922457419601
$ws.Range("E6").Value = 0.7874026932376078
$ws.Range("F6").Value = 0.8262992831541218
$ws.Range("G6").Value = 0.8167783491766661
$ws.Range("H6").Value = 0.8861239446690317
$ws.Range("I6").Value = 0.8197978904605989
$ws.Range("J6").Value = 0.02807391593598027
$ws.Range("L6").Value = 0.8074029465158497
$ws.Range("M6").Value = 0.8185147007303205
$ws.Range("N6").Value = 0.8226325905147783
$ws.Range("O6").Value = 0.8683393070489845
$ws.Range("P6").Value = 0.7630757257173136
$ws.Range("Q6").Value = 0.0384425666860263
$ws.Range("R6").Value = 0.7630446345975949
$ws.Range("S6").Value = 0.738536368830899
$ws.Range("T6").Value = 0.7967229902713775
$ws.Range("U6").Value = 0.7055578829772378
$ws.Range("V6").Value = 0.8115167519094587
$ws.Range("W6").Value = 0.8115330914249327
$ws.Range("X6").Value = 0.01693252094512161
$ws.Range("Y6").Value = 0.7920449863571701
$ws.Range("Z6").Value = 0.8047962342317182
$ws.Range("AA6").Value = 0.8156082467910425
$ws.Range("AB6").Value = 0.8032887560181426
$ws.Range("AC6").Value = 0.8419272337265902
$ws.Range("AD6").Value = 0.8063140172651634
$ws.Range("AE6").Value = 0.0192611199657048
$ws.Range("AF6").Value = 0.8007920906567992
$ws.Range("AG6").Value = 0.7967701205604431
$ws.Range("AH6").Value = 0.8192683770357795
$ws.Range("AI6").Value = 0.7794114631655178
$ws.Range("AJ6").Value = 0.8353280349072776
$ws.Range("AK6").Value = 0.8217103988245121
$ws.Range("AL6").Value = 0.04990396097340202
$ws.Range("AM6").Value = 0.784726443768997
$ws.Range("AN6").Value = 0.7588538999829322
$ws.Range("AO6").Value = 0.862072744585478
$ws.Range("AP6").Value = 0.8082329479613012
$ws.Range("AQ6").Value = 0.8946659578238525
$ws.Range("B7").Value = 0.8580958377807846
$ws.Range("C7").Value = 0.04104877307437271
$ws.Range("D7").Value = 0.8463916884409605
$ws.Range("F7").Value = 0.8958094160131512
$ws.Range("H7").Value = 0.9158724340175952
$ws.Range("I7").Value = 0.8755042041822045
$ws.Range("J7").Value = 0.04413844179317703
$ws.Range("K7").Value = 0.8626757496578326
$ws.Range("L7").Value = 0.8208922844746274
$ws.Range("M7").Value = 0.91711357246162
$ws.Range("N7").Value = 0.8408328432297273
$ws.Range("P7").Value = 0.8425791006061886
$ws.Range("Q7").Value = 0.03757084580740404
$ws.Range("R7").Value = 0.8279323513366067
$ws.Range("S7").Value = 0.8170858943978888
$ws.Range("T7").Value = 0.829047192353644
$ws.Range("W7").Value = 0.8845616733037713
$ws.Range("X7").Value = 0.05169178699104997
$ws.Range("Y7").Value = 0.8904086457277947
$ws.Range("Z7").Value = 0.7936049801924165
$ws.Range("AC7").Value = 0.9272840479334539
$ws.Range("AD7").Value = 0.8659991445317573
$ws.Range("AE7").Value = 0.03365642527193489
$ws.Range("AG7").Value = 0.8175243215565795
$ws.Range("AJ7").Value = 0.8835125448028672
$ws.Range("AK7").Value = 0.8707042560563728
$ws.Range("AL7").Value = 0.04696120126963177
$ws.Range("AM7").Value = 0.898923727759114
$ws.Range("AN7").Value = 0.7836082933456348
$ws.Range("AP7").Value = 0.8642058974127286
